# Slide 1, shape 2 ("文字方塊 5" / Text Box 5, id=6) is resized/repositioned
# and its run bumped from 44pt to 54pt.
#
# PowerPoint's Shape.Left/Top/Width/Height are expressed in points, but the
# underlying COM host stores them as 32-bit floats, and the point->EMU
# conversion truncates rather than rounds. The literals below are the
# closest float32-representable point values that round-trip to the exact
# target EMU (444569, 267279, 1260459, 923330) once converted, so the
# resulting OOXML matches byte-for-byte rather than being off by a
# stray EMU.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item(2)

$shape.Left = 35.0054359436
$shape.Top = 21.0455913544
$shape.Width = 99.2487411499
$shape.Height = 72.7031555176

$shape.TextFrame.TextRange.Font.Size = 54
